$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 80 (shifts old row 80 -> row 81, keeps formatting)
$ws.Rows.Item(80).Insert()

# Fill in the new row 80 with data (same market/company/region as row 81, new date & prices)
$ws.Range("A80").Value = 1
$ws.Range("B80").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C80").Value = "Arica y Parinacota"
$ws.Range("D80").Value = 45239
$ws.Range("E80").Value = 15
$ws.Range("F80").Value = 100112052
$ws.Range("G80").Value = "Albahaca"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 370
$ws.Range("K80").Value = 1300
$ws.Range("L80").Value = 1500
$ws.Range("M80").Value = 1381
$ws.Range("N80").Value = "$/paquete"
$ws.Range("O80").Value = "Región de Arica y Parinacota"
$ws.Range("P80").Value = 1381
$ws.Range("Q80").Value = 1
$ws.Range("R80").Value = "Hortaliza"
